# Qr View for Workflow
# Adds a new permission/resource row (row 9) to the DATA-PERM sheet,
# mirroring the pattern used by existing rows 6-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA-PERM")

# New row 9: same permId as rows 6-8, new resourceId, new name (查询视图读取)
$ws.Range("A9").Value = "e501b47a-c08b-4c83-b12b-95ad82873e96"
$ws.Range("B9").Value = "2058d59f-a2d1-43ee-b6c1-9687d6018f61"
$ws.Range("C9").Value = "查询视图读取"

# Copy the styling from the row above (row 8) so row 9 matches the
# formatting used by the other permission rows.
$ws.Range("A8:C8").Copy() | Out-Null
$ws.Range("A9:C9").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

# Re-set the values (paste of formats only should not touch them, but just
# to be safe re-apply after the formatting paste)
$ws.Range("A9").Value = "e501b47a-c08b-4c83-b12b-95ad82873e96"
$ws.Range("B9").Value = "2058d59f-a2d1-43ee-b6c1-9687d6018f61"
$ws.Range("C9").Value = "查询视图读取"

# Update the active selection to match the saved cursor position (C10)
$ws.Range("C10").Select() | Out-Null

$wb.Save()
